$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Lab 7 (Branch&Bound) grade for the student in row 6
$ws.Range("H6").Value = 6

# Feedback comment for Lab 7 column (row 7, merged H7:H14)
$ws.Range("H7").Value = "Good code but PDF document not submitted"

# Test mark (column I) becomes "Passed" instead of blank
$ws.Range("I6").Value = "Passed"

# Update the active selection to reflect the next grading cell (I7:I14)
$ws.Range("I7:I14").Select()
